$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "C1"
$ws.Range("B2").Value = '${B1}'
$ws.Range("C2").Value = '${C1_1}-${C1_2}'

$ws.Range("C2").Select()
